$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuille1")

$ws.Range("C34").Value = 34953667
$ws.Range("D34").Value = 1369862
$ws.Range("E34").Value = 224360
$ws.Range("F34").Value = 155
$ws.Range("G34").Value = "IDF*"

$ws.Range("G51:H51").Copy()
$ws.Range("G52:H52").PasteSpecial(-4122)
$ws.Range("G52").Value = "IDF*"
$ws.Range("H52").Value = "Ile-de-France sur PC propre"
